$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leave Application - Hourly: update TestCases value in row 2 from 48 to 42
$ws.Range("B2").Value = "42"

# Update the active selection to C4
$ws.Range("C4").Select() | Out-Null
